{"js": "// Update two-digit multiplication problems in the worksheet table.\nconst replacements = [\n  [\"89\u00d717=\", \"35\u00d711=\"],\n  [\"94\u00d769=\", \"94\u00d748=\"],\n  [\"63\u00d731=\", \"47\u00d721=\"],\n  [\"35\u00d751=\", \"58\u00d780=\"],\n  [\"18\u00d744=\", \"29\u00d797=\"],\n  [\"13\u00d765=\", \"92\u00d758=\"],\n  [\"54\u00d740=\", \"75\u00d743=\"],\n  [\"13\u00d744=\", \"50\u00d761=\"],\n  [\"67\u00d765=\", \"80\u00d799=\"],\n  [\"20\u00d731=\", \"51\u00d779=\"],\n  [\"32\u00d777=\", \"96\u00d733=\"],\n  [\"43\u00d773=\", \"64\u00d789=\"],\n  [\"84\u00d790=\", \"57\u00d762=\"],\n  [\"78\u00d764=\", \"84\u00d775=\"],\n  [\"12\u00d757=\", \"85\u00d773=\"],\n  [\"46\u00d729=\", \"54\u00d730=\"],\n  [\"40\u00d748=\", \"96\u00d752=\"],\n  [\"13\u00d750=\", \"41\u00d721=\"],\n  [\"46\u00d716=\", \"53\u00d787=\"],\n  [\"79\u00d756=\", \"36\u00d731=\"],\n  [\"46\u00d788=\", \"39\u00d758=\"],\n  [\"46\u00d713=\", \"88\u00d742=\"],\n  [\"37\u00d768=\", \"91\u00d755=\"],\n  [\"73\u00d752=\", \"63\u00d786=\"],\n  [\"53\u00d781=\", \"26\u00d733=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "# Update two-digit multiplication problems in the worksheet table.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"89\u00d717=\", \"35\u00d711=\"),\n    @(\"94\u00d769=\", \"94\u00d748=\"),\n    @(\"63\u00d731=\", \"47\u00d721=\"),\n    @(\"35\u00d751=\", \"58\u00d780=\"),\n    @(\"18\u00d744=\", \"29\u00d797=\"),\n    @(\"13\u00d765=\", \"92\u00d758=\"),\n    @(\"54\u00d740=\", \"75\u00d743=\"),\n    @(\"13\u00d744=\", \"50\u00d761=\"),\n    @(\"67\u00d765=\", \"80\u00d799=\"),\n    @(\"20\u00d731=\", \"51\u00d779=\"),\n    @(\"32\u00d777=\", \"96\u00d733=\"),\n    @(\"43\u00d773=\", \"64\u00d789=\"),\n    @(\"84\u00d790=\", \"57\u00d762=\"),\n    @(\"78\u00d764=\", \"84\u00d775=\"),\n    @(\"12\u00d757=\", \"85\u00d773=\"),\n    @(\"46\u00d729=\", \"54\u00d730=\"),\n    @(\"40\u00d748=\", \"96\u00d752=\"),\n    @(\"13\u00d750=\", \"41\u00d721=\"),\n    @(\"46\u00d716=\", \"53\u00d787=\"),\n    @(\"79\u00d756=\", \"36\u00d731=\"),\n    @(\"46\u00d788=\", \"39\u00d758=\"),\n    @(\"46\u00d713=\", \"88\u00d742=\"),\n    @(\"37\u00d768=\", \"91\u00d755=\"),\n    @(\"73\u00d752=\", \"63\u00d786=\"),\n    @(\"53\u00d781=\", \"26\u00d733=\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText, $false, $false, $false, $false, $false,\n        $true, 1, $true, $newText, 2\n    ) | Out-Null\n}"}
